$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# ---------------------------------------------------------------------------
# 1) Header row: shrink the italic header labels from 11pt (sz 22) to
#    10pt (sz 20). Apply to each header cell's own text range (excluding the
#    trailing cell/paragraph mark) so the paragraph-mark run properties are
#    left untouched.
# ---------------------------------------------------------------------------
for ($c = 1; $c -le $t.Columns.Count; $c++) {
    $cell = $t.Cell(1, $c)
    $cr = $cell.Range
    $textRange = $d.Range($cr.Start, $cr.End - 1)
    $textRange.Font.Size = 10
}

# ---------------------------------------------------------------------------
# 2) Data rows: right-align the stock-name (first) column instead of left.
# ---------------------------------------------------------------------------
for ($r = 2; $r -le $t.Rows.Count; $r++) {
    $cell = $t.Cell($r, 1)
    $cell.Range.ParagraphFormat.Alignment = 2
}

# ---------------------------------------------------------------------------
# 3) Reformat numeric values to three decimal places in specific cells.
# ---------------------------------------------------------------------------
function Set-CellValue($table, $row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $cr = $cell.Range
    $textRange = $d.Range($cr.Start, $cr.End - 1)
    $textRange.Text = $newText
}

Set-CellValue $t 2 2 "1.630"
Set-CellValue $t 2 3 "3.000"
Set-CellValue $t 2 5 "1.000"

Set-CellValue $t 3 2 "0.350"
Set-CellValue $t 3 3 "0.000"
Set-CellValue $t 3 4 "1.000"
Set-CellValue $t 3 5 "0.000"
Set-CellValue $t 3 6 "0.500"

Set-CellValue $t 5 3 "3.000"
Set-CellValue $t 5 5 "1.000"
Set-CellValue $t 5 6 "1.000"

Set-CellValue $t 6 3 "2.000"

Set-CellValue $t 8 3 "0.000"
Set-CellValue $t 8 5 "0.000"

Set-CellValue $t 9 3 "2.000"
Set-CellValue $t 9 6 "0.810"

Set-CellValue $t 10 2 "2.980"
Set-CellValue $t 10 3 "0.000"
Set-CellValue $t 10 5 "0.000"

Set-CellValue $t 11 2 "5.010"
Set-CellValue $t 11 3 "0.000"
Set-CellValue $t 11 4 "0.000"
Set-CellValue $t 11 5 "0.000"
Set-CellValue $t 11 6 "0.000"

Set-CellValue $t 12 2 "3.560"
Set-CellValue $t 12 3 "0.660"
Set-CellValue $t 12 5 "0.220"
